$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update jornada labels JOR. 6 -> JOR. 7
$ws.Range("E2").Value = "JOR. 7"
$ws.Range("E3").Value = "JOR. 7"
$ws.Range("E4").Value = "JOR. 7"

# Update matchups for row 2
$ws.Range("G2").Value = "Armada"
$ws.Range("H2").Value = "Puche"

# Update matchups for row 3
$ws.Range("G3").Value = "Gonzo"
$ws.Range("H3").Value = "Coquina"

# Update matchups for row 4 (G4 stays "Papu")
$ws.Range("H4").Value = "Kike"
